$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '65.040.27'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '3.156.09'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.76'
$ws.Range("E5").Value = '  +1.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.56'
$ws.Range("D8").Value = '3.155.92'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("E10").Value = '  -1.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.14'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.500'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.15'
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").Value = '3.677.84'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").Value = '64.958.36'
$ws.Range("D17").Value = '3.164.92'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.14'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '504.69'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.91'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.715'
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.21'
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.74'
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.48'
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.07'
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("E30").Value = '  +4.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.57'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.43'
$ws.Range("E33").Value = '  +4.36%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.49'
$ws.Range("E35").Value = '  -1.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.06'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0887'
$ws.Range("E37").Value = '  +3.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '478.14'
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0415'
$ws.Range("E39").Value = '  -1.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.93'
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.73'
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("D42").Value = '2.991.20'
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("E43").Value = '  -1.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.282'
$ws.Range("E44").Value = '  -3.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  -1.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.28'
$ws.Range("E46").Value = '  -4.08%  '
$ws.Range("D47").Value = '0.0₃0592'
$ws.Range("E47").Value = '  +3.38%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("E51").Value = '  +14.66%  '
